$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds student codes that must stay TEXT (several look like plain
# numbers, and a couple - "0566", "01801" - have leading zeros that would be
# silently dropped if Excel stored them as numbers). For every row, force
# the cell to Text before writing the new value, then restore the default
# ("Normal") cell style afterwards so formatting ends up matching the rest
# of the sheet (only the cell's stored content/type should change, per the
# diff).
$codes = @{
  2  = "31986"
  3  = "113333"
  4  = "1111768"
  5  = "10865"
  6  = "11274"
  7  = "0566"
  8  = "01801"
  9  = "18000600"
  10 = "80056"
  11 = "20027"
  12 = "101552"
  13 = "11207"
  14 = "10005"
  15 = "812"
  16 = "1165"
  17 = "17023"
  18 = "8101172"
}

foreach ($row in 2..18) {
  $cell = $ws.Cells.Item($row, 1)
  $cell.NumberFormat = "@"
  $cell.Value = $codes[$row]
  $cell.Style = "Normal"
}

# Other numeric score updates
$ws.Range("D5").Value  = 3
$ws.Range("B10").Value = 2
$ws.Range("C13").Value = 3
$ws.Range("D15").Value = 3
